$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 97, shifting existing rows 97-151 down to 98-152
$ws.Rows(97).Insert()

# Populate the newly inserted row 97 with the new record's data
$ws.Range("A97").Value = 5
$ws.Range("B97").Value = "Macroferia Regional de Talca"
$ws.Range("C97").Value = "Maule"
$ws.Range("D97").Value = 45236
$ws.Range("E97").Value = 7
$ws.Range("F97").Value = 100112022
$ws.Range("G97").Value = "Arveja Verde"
$ws.Range("H97").Value = "Sin especificar"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 400
$ws.Range("K97").Value = 20000
$ws.Range("L97").Value = 23000
$ws.Range("M97").Value = 21500
$ws.Range("N97").Value = "$/saco 25 kilos"
$ws.Range("O97").Value = "Región del Maule"
$ws.Range("P97").Value = 860
$ws.Range("Q97").Value = 25
$ws.Range("R97").Value = "Hortaliza"
